# Admin UI TC implemented
# Adds a new "campaignManagementTest" test case row to both the
# RUNMANAGER and DATA sheets, and flips the "execute" flag for the
# other login-related test rows from "yes" to "no".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RUNMANAGER
$ws2 = $wb.Worksheets.Item(2)   # DATA

# ---------------------------------------------------------------
# RUNMANAGER sheet
# ---------------------------------------------------------------

# loginlogoutTest is no longer executed
$ws1.Range("C2").Value = "no"

# Append a new row for the campaign management test case.
# The priority/count columns hold numeric-looking text ("1") entered
# with a leading quote prefix, matching the rest of the column.
$ws1.Range("A7").Value = "campaignManagementTest"
$ws1.Range("C7").Value = "yes"
$ws1.Range("D7").Value = "'1"
$ws1.Range("E7").Value = "'1"

# ---------------------------------------------------------------
# DATA sheet
# ---------------------------------------------------------------

# These rows are no longer executed
$ws2.Range("B2").Value = "no"
$ws2.Range("B5").Value = "no"
$ws2.Range("B8").Value = "no"

# Append a new data row for the campaign management test case.
$ws2.Range("A9").Value = "campaignManagementTest"
$ws2.Range("B9").Value = "yes"
$ws2.Range("C9").Value = "chrome"
$ws2.Range("D9").Value = "'96"
$ws2.Range("E9").Value = "'"
$ws2.Range("F9").Value = "'"
$ws2.Range("G9").Value = "'"

# Back to RUNMANAGER for the test description (keeps shared-string
# insertion order aligned with the saved workbook).
$ws1.Range("B7").Value = "To check if the user is able to create a campaign"

# ---------------------------------------------------------------
# Selection / active sheet bookkeeping (matches the saved view state)
# ---------------------------------------------------------------

$ws1.Activate()
$ws1.Range("I13").Select()

$ws2.Activate()
$ws2.Range("D9").Select()
